$d = $word.ActiveDocument

# 1. Update the DATE field's cached result text (July 15 -> July 16, 2015)
$d.Content.Find.Execute("July 15, 2015", $true, $false, $false, $false, $false, $true, 1, $false, "July 16, 2015", 2) | Out-Null

# 2. Rewrite the opening paragraph of the Introduction
$old1 = "Many experimental investigators in microbiology, genetics, and evolutionary biology use growth curves to estimate fitness. They measure the Optical Density (OD) of one or more populations of cells over several hours or even days to acquire the growth curves. The simplest way to estimate fitness from these curves is to infer the growth rate: taking the log of the curves during the exponential growth phase, using linear regression to fit a linear line to the data, and taking the slope of the line as a measure of the growth rate"
$new1 = "Growth curves are common method to estimate fitness in microbiology, genetics, and evolutionary biology. Optical Density (OD) is used to acquire the growth curves of one or more populations of cells, over varying time periods. The simplest way to estimate fitness from these curves is to infer the growth rate during the exponential growth phase. This is done by taking the log of the mean of the growth curves during the exponential growth phase and using linear regression to estimate the slope of the curves as a measure of the growth rate"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 3. Note that the additional growth phases affect the selection coefficient
$old2 = "growth phases in addition to"
$new2 = "growth phases that affect the selection coefficient in addition to"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 4. Rewrite the competition-assays paragraph
$old3 = "Competition assays are a common fitness inference method that takes these additional growth phases into account. Competition assays include the growth of two strains in the same container – the strain of interest and a reference strain (for example, a mutant strain and a wildtype strain). From the change in frequency over the competition"
$new3 = "Pairwise competition experiments are commonly used to infer fitness in a way that takes these additional growth phases into account. Competition experiments include the growth of two strains in the same container – the strain of interest and a reference strain (for example, a mutant strain and a wildtype strain). The frequency of each strain is measured during the experiment. From the change in frequencies over the competition"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 5. "Because competition assays require" -> "Because pairwise competition experiments require"
$old4 = "Because competition assays require"
$new4 = "Because pairwise competition experiments require"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# 6. Drop "can be used to" and relocate the _GoBack bookmark to that spot
$old5 = "for population genetics models that can be used to test hypotheses"
$new5 = "for population genetics models that test hypotheses"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

$bmRange = $d.Content
$bmRange.Find.Execute("for population genetics models that ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
